# Updated cryptos list with GitHub Actions - applies latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text assignments (E column percentage strings and non-numeric-looking D values)
$ws.Range("D2").Value = '70.729.87'
$ws.Range("E2").Value = '  -2.18%  '
$ws.Range("D3").Value = '3.634.65'
$ws.Range("E3").Value = '  +0.24%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("E5").Value = '  -2.10%  '
$ws.Range("E6").Value = '  -3.44%  '
$ws.Range("D7").Value = '3.625.94'
$ws.Range("E7").Value = '  +0.29%  '
$ws.Range("E8").Value = '  +1.77%  '
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("E10").Value = '  -4.88%  '
$ws.Range("E11").Value = '  +17.92%  '
$ws.Range("E12").Value = '  +0.20%  '
$ws.Range("E13").Value = '  -3.58%  '
$ws.Range("E14").Value = '  -1.99%  '
$ws.Range("D15").Value = '4.220.33'
$ws.Range("E15").Value = '  +0.24%  '
$ws.Range("E16").Value = '  -4.54%  '
$ws.Range("E17").Value = '  +0.22%  '
$ws.Range("D18").Value = '3.635.13'
$ws.Range("E18").Value = '  +0.00%  '
$ws.Range("D19").Value = '70.789.39'
$ws.Range("E19").Value = '  -2.23%  '
$ws.Range("E20").Value = '  -0.37%  '
$ws.Range("E21").Value = '  -4.22%  '
$ws.Range("E22").Value = '  -1.75%  '
$ws.Range("E23").Value = '  +0.98%  '
$ws.Range("E24").Value = '  -3.53%  '
$ws.Range("E25").Value = '  -4.44%  '
$ws.Range("E26").Value = '  -2.64%  '
$ws.Range("E27").Value = '  -2.68%  '
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("E29").Value = '  -2.09%  '
$ws.Range("E30").Value = '  -2.34%  '
$ws.Range("E31").Value = '  +0.54%  '
$ws.Range("E32").Value = '  -5.82%  '
$ws.Range("E33").Value = '  +0.93%  '
$ws.Range("E34").Value = '  -6.54%  '
$ws.Range("E35").Value = '  -4.50%  '
$ws.Range("E36").Value = '  -3.26%  '
$ws.Range("E37").Value = '  -1.90%  '
$ws.Range("E38").Value = '  -0.75%  '
$ws.Range("E39").Value = '  -2.18%  '
$ws.Range("E40").Value = '  +0.06%  '
$ws.Range("E41").Value = '  +0.59%  '
$ws.Range("D42").Value = '3.546.55'
$ws.Range("E42").Value = '  -2.40%  '
$ws.Range("E43").Value = '  -1.66%  '
$ws.Range("E44").Value = '  -3.58%  '
$ws.Range("E45").Value = '  -4.03%  '
$ws.Range("E46").Value = '  -5.52%  '
$ws.Range("E49").Value = '  +1.27%  '
$ws.Range("E50").Value = '  +2.94%  '
$ws.Range("E51").Value = '  -3.35%  '

# D column values that look like pure numbers must be forced to remain text
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '584.01'
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '176.04'
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.618'
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '6.82'
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.609'
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '48.54'
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.0000284'
$c.Style = "Normal"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '674.49'
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '8.99'
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '17.80'
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '11.50'
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.943'
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '17.20'
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '99.91'
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '2.79'
$c.Style = "Normal"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '34.62'
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '9.17'
$c.Style = "Normal"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '7.54'
$c.Style = "Normal"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '3.97'
$c.Style = "Normal"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '575.80'
$c.Style = "Normal"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '11.10'
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '58.51'
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '34.41'
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '137.75'
$c.Style = "Normal"

# Rows 47 and 48 swapped positions (ThetaToken now ranks above Fetch.AI)
$ws.Range("B47").Value = 'ThetaToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '2.96'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +5.12%  '

$ws.Range("B48").Value = 'Fetch.AI'
$ws.Range("C48").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '2.68'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -4.40%  '
